$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 94
$ws.Range("I4").Value = 95
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 95
$ws.Range("L4").Value = 90
$ws.Range("M4").Value = 19
$ws.Range("N4").Value = -318
$ws.Range("H13").Value = 1000000
$ws.Range("I13").Value = 1000000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1000000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -999831
$ws.Range("N13").ClearContents()
$ws.Range("H15").Value = 248.53334
$ws.Range("I15").Value = 248.53334
$ws.Range("K15").Value = 745.6000200000001
$ws.Range("M15").Value = -576.6000200000001
$ws.Range("H32").Value = 2722.4666
$ws.Range("I32").Value = 3590.6667
$ws.Range("J32").Value = 2143.6667
$ws.Range("K32").Value = 3590.6667
$ws.Range("L32").Value = 2143.6667
$ws.Range("M32").Value = -3264.6667
$ws.Range("N32").Value = -2795.6667
$ws.Range("H69").Value = 10206.5
$ws.Range("I69").Value = 413
$ws.Range("J69").Value = 20000
$ws.Range("K69").Value = 1239
$ws.Range("L69").Value = 60000
$ws.Range("M69").Value = -365
$ws.Range("N69").Value = -61748
$ws.Range("H72").Value = 10206.5
$ws.Range("I72").Value = 413
$ws.Range("J72").Value = 20000
$ws.Range("K72").Value = 3717
$ws.Range("L72").Value = 180000
$ws.Range("M72").Value = 651
$ws.Range("N72").Value = -188736
$ws.Range("H112").Value = 35699.965
$ws.Range("J112").Value = 39518.703
$ws.Range("L112").Value = 118556.109
$ws.Range("N112").Value = -120772.109
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8487.134
$ws.Range("I32").Value = 5543.024
$ws.Range("K32").Value = 5543.024
$ws.Range("M32").Value = -5256.024
$ws.Range("H74").Value = 58345.24
$ws.Range("I74").Value = 79027.30499999999
$ws.Range("K74").Value = 79027.30499999999
$ws.Range("M74").Value = -78153.30499999999
$ws.Range("H77").Value = 58345.24
$ws.Range("I77").Value = 79027.30499999999
$ws.Range("K77").Value = 395136.525
$ws.Range("M77").Value = -390768.525
$ws.Range("H92").Value = 49999
$ws.Range("J92").Value = 49999
$ws.Range("L92").Value = 49999
$ws.Range("N92").Value = -54991
$ws.Range("H110").Value = 8522.814
$ws.Range("I110").Value = 8664.639999999999
$ws.Range("K110").Value = 8664.639999999999
$ws.Range("M110").Value = -6619.639999999999
$ws.Range("H122").Value = 1133.579
$ws.Range("I122").Value = 867.4666999999999
$ws.Range("K122").Value = 2602.4001
$ws.Range("M122").Value = -152.4000999999998
$ws.Range("H132").Value = 4231.2173
$ws.Range("I132").Value = 4062.762
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 12188.286
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -9658.286
$ws.Range("N132").Value = -23060
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 16814.166
$ws.Range("J106").Value = 16814.166
$ws.Range("L106").Value = 16814.166
$ws.Range("N106").Value = -19338.166
$ws.Range("H107").Value = 1141.1818
$ws.Range("I107").Value = 1100.375
$ws.Range("K107").Value = 1100.375
$ws.Range("M107").Value = 819.625
$ws.Range("H134").Value = 1907.7391
$ws.Range("I134").Value = 1906.4777
$ws.Range("K134").Value = 5719.4331
$ws.Range("M134").Value = -3184.4331
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2666.2654
$ws.Range("I132").Value = 2579.2827
$ws.Range("K132").Value = 7737.848100000001
$ws.Range("M132").Value = -5207.848100000001
$ws.Range("H134").Value = 11670.925
$ws.Range("I134").Value = 4450.2446
$ws.Range("K134").Value = 13350.7338
$ws.Range("M134").Value = -10815.7338
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 230.71428
$ws.Range("I11").Value = 213
$ws.Range("J11").Value = 275
$ws.Range("K11").Value = 639
$ws.Range("L11").Value = 825
$ws.Range("M11").Value = -499
$ws.Range("N11").Value = -1105
$ws.Range("H80").Value = 3
$ws.Range("J80").Value = 3
$ws.Range("L80").Value = 9
$ws.Range("N80").Value = -1881
$ws.Range("H83").Value = 3
$ws.Range("J83").Value = 3
$ws.Range("L83").Value = 27
$ws.Range("N83").Value = -9387
$ws.Range("H86").Value = 963.1667
$ws.Range("J86").Value = 1269.75
$ws.Range("L86").Value = 3809.25
$ws.Range("N86").Value = -6181.25
$ws.Range("H87").Value = 23436.637
$ws.Range("I87").Value = 19686.285
$ws.Range("J87").Value = 29999.75
$ws.Range("K87").Value = 59058.855
$ws.Range("L87").Value = 89999.25
$ws.Range("M87").Value = -57810.855
$ws.Range("N87").Value = -92495.25
$ws.Range("H89").Value = 963.1667
$ws.Range("J89").Value = 1269.75
$ws.Range("L89").Value = 11427.75
$ws.Range("N89").Value = -23283.75
$ws.Range("H90").Value = 23436.637
$ws.Range("I90").Value = 19686.285
$ws.Range("J90").Value = 29999.75
$ws.Range("K90").Value = 177176.565
$ws.Range("L90").Value = 269997.75
$ws.Range("M90").Value = -170936.565
$ws.Range("N90").Value = -282477.75
$ws.Range("H113").Value = 866.0769
$ws.Range("I113").Value = 783.3333
$ws.Range("K113").Value = 2349.9999
$ws.Range("M113").Value = -179.9998999999998
$ws.Range("H115").Value = 28
$ws.Range("I115").Value = 28
$ws.Range("K115").Value = 84
$ws.Range("M115").Value = 1091
$ws.Range("H131").Value = 34225.805
$ws.Range("I131").Value = 251250
$ws.Range("J131").Value = 2074.074
$ws.Range("K131").Value = 753750
$ws.Range("L131").Value = 6222.222
$ws.Range("M131").Value = -748710
$ws.Range("N131").Value = -16302.222
$ws.Range("H132").Value = 1374.3889
$ws.Range("J132").Value = 1399.4166
$ws.Range("L132").Value = 12594.7494
$ws.Range("N132").Value = -17654.7494
$ws.Range("H140").Value = 4607.8335
$ws.Range("I140").Value = 4661.75
$ws.Range("K140").Value = 13985.25
$ws.Range("M140").Value = -8805.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 37063696
$ws.Range("I18").Value = 111111110
$ws.Range("K18").Value = 111111110
$ws.Range("M18").Value = -111110817
$ws.Range("H102").Value = 58826540
$ws.Range("I102").Value = 3038.5
$ws.Range("K102").Value = 3038.5
$ws.Range("M102").Value = -1416.5
$ws.Range("H122").Value = 2610
$ws.Range("J122").Value = 4229.8335
$ws.Range("L122").Value = 12689.5005
$ws.Range("N122").Value = -17589.5005
$ws.Range("H123").Value = 40162.5
$ws.Range("J123").Value = 40162.5
$ws.Range("L123").Value = 40162.5
$ws.Range("N123").Value = -45062.5
$ws.Range("H126").Value = 15656
$ws.Range("I126").Value = 19281.54
$ws.Range("K126").Value = 57844.62
$ws.Range("M126").Value = -55374.62
$ws.Range("H132").Value = 3190.1
$ws.Range("I132").Value = 2492.6
$ws.Range("J132").Value = 5282.6
$ws.Range("K132").Value = 7477.799999999999
$ws.Range("L132").Value = 15847.8
$ws.Range("M132").Value = -4947.799999999999
$ws.Range("N132").Value = -20907.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 7113.5835
$ws.Range("I13").Value = 6920
$ws.Range("K13").Value = 6920
$ws.Range("M13").Value = -6780
$ws.Range("H22").Value = 2949.4
$ws.Range("I22").Value = 2984.8572
$ws.Range("J22").Value = 2866.6667
$ws.Range("K22").Value = 2984.8572
$ws.Range("L22").Value = 2866.6667
$ws.Range("M22").Value = -2689.8572
$ws.Range("N22").Value = -3456.6667
$ws.Range("H27").Value = 2949.4
$ws.Range("I27").Value = 2984.8572
$ws.Range("J27").Value = 2866.6667
$ws.Range("K27").Value = 2984.8572
$ws.Range("L27").Value = 2866.6667
$ws.Range("M27").Value = -2877.8572
$ws.Range("N27").Value = -3080.6667
$ws.Range("H61").Value = 3616.7083
$ws.Range("I61").Value = 3609.1428
$ws.Range("J61").Value = 3669.6667
$ws.Range("K61").Value = 3609.1428
$ws.Range("L61").Value = 3669.6667
$ws.Range("M61").Value = -3407.1428
$ws.Range("N61").Value = -4073.6667
$ws.Range("H113").Value = 3616.7083
$ws.Range("I113").Value = 3609.1428
$ws.Range("J113").Value = 3669.6667
$ws.Range("K113").Value = 3609.1428
$ws.Range("L113").Value = 3669.6667
$ws.Range("M113").Value = -1439.1428
$ws.Range("N113").Value = -8009.6667
$ws.Range("H122").Value = 261600.48
$ws.Range("I122").Value = 338054
$ws.Range("J122").Value = 6755.4443
$ws.Range("K122").Value = 1014162
$ws.Range("L122").Value = 20266.3329
$ws.Range("M122").Value = -1011712
$ws.Range("N122").Value = -25166.3329
$ws.Range("H132").Value = 4827.1816
$ws.Range("I132").Value = 4409.9
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 13229.7
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -10699.7
$ws.Range("N132").Value = -32060
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 70000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H122").Value = 2540.8823
$ws.Range("I122").Value = 2168.8462
$ws.Range("J122").Value = 3750
$ws.Range("K122").Value = 6506.5386
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -4056.5386
$ws.Range("N122").Value = -16150
$ws.Range("H132").Value = 2868.3333
$ws.Range("I132").Value = 2868.3333
$ws.Range("K132").Value = 8604.999899999999
$ws.Range("M132").Value = -6074.999899999999
$ws.Range("H136").Value = 2796.05
$ws.Range("I136").Value = 2894.4666
$ws.Range("J136").Value = 2500.8
$ws.Range("K136").Value = 8683.399800000001
$ws.Range("L136").Value = 7502.400000000001
$ws.Range("M136").Value = -6133.399800000001
$ws.Range("N136").Value = -12602.4
